$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1) to reflect the Nov 2020 reporting period
$ws.Range("C1").Value = "AN_ANUAL_202011"
$ws.Range("D1").Value = "AN_SEM_202011"
$ws.Range("E1").Value = "AN_TRI_202011"
$ws.Range("F1").Value = "AC_20201115"
$ws.Range("G1").Value = "AN_20201115"
$ws.Range("L1").Value = "AN_COTA_202009"
$ws.Range("M1").Value = "AN_COTA_202010"
$ws.Range("N1").Value = "AN_COTA_20201115"

# Update measurement data (rows 2-42) for columns C-G (monitoring indicators)
# and L-N (AN_COTA indicators), per the 16 Nov 2020 data refresh
$ws.Range("C2").Value = -27.4
$ws.Range("D2").Value = -38.7
$ws.Range("E2").Value = -76.3
$ws.Range("F2").Value = 4.7
$ws.Range("G2").Value = -97
$ws.Range("L2").Value = -5.68
$ws.Range("M2").Value = -9.220000000000001
$ws.Range("N2").Value = -11.55
$ws.Range("C3").Value = -41
$ws.Range("D3").Value = -25.3
$ws.Range("E3").Value = -53.4
$ws.Range("F3").Value = 53.5
$ws.Range("G3").Value = -55.1
$ws.Range("L3").Value = -54.38
$ws.Range("M3").Value = -51.9
$ws.Range("N3").Value = -31.19
$ws.Range("C4").Value = -40.8
$ws.Range("D4").Value = -46.5
$ws.Range("E4").Value = -84
$ws.Range("F4").Value = 34.5
$ws.Range("G4").Value = -83
$ws.Range("C5").Value = -18.8
$ws.Range("D5").Value = -35.3
$ws.Range("E5").Value = -70.7
$ws.Range("F5").Value = 24.4
$ws.Range("G5").Value = -79.8
$ws.Range("L5").Value = -46.96
$ws.Range("M5").Value = -64.04000000000001
$ws.Range("N5").Value = -59.66
$ws.Range("C6").Value = -34
$ws.Range("D6").Value = -60
$ws.Range("E6").Value = -87
$ws.Range("F6").Value = 13.6
$ws.Range("G6").Value = -88.90000000000001
$ws.Range("L6").Value = -6.29
$ws.Range("M6").Value = -4.1
$ws.Range("N6").Value = -16.53
$ws.Range("C7").Value = -27.6
$ws.Range("D7").Value = -42.7
$ws.Range("E7").Value = -71.59999999999999
$ws.Range("F7").Value = 32.4
$ws.Range("G7").Value = -76.3
$ws.Range("C8").Value = -52.9
$ws.Range("D8").Value = -35.1
$ws.Range("E8").Value = -64.40000000000001
$ws.Range("F8").Value = 65.59999999999999
$ws.Range("G8").Value = -47.1
$ws.Range("L8").Value = -5.66
$ws.Range("M8").Value = -4.55
$ws.Range("N8").Value = 10.46
$ws.Range("C9").Value = -33.5
$ws.Range("D9").Value = -46
$ws.Range("E9").Value = -71.40000000000001
$ws.Range("F9").Value = 27.9
$ws.Range("G9").Value = -82.5
$ws.Range("C10").Value = -52.4
$ws.Range("D10").Value = -41
$ws.Range("E10").Value = -56.5
$ws.Range("F10").Value = 58.1
$ws.Range("G10").Value = -46
$ws.Range("L10").Value = -22.57
$ws.Range("M10").Value = -42.3
$ws.Range("N10").Value = -33.04
$ws.Range("C11").Value = -48.9
$ws.Range("D11").Value = -35.1
$ws.Range("E11").Value = -50.8
$ws.Range("F11").Value = 87.8
$ws.Range("G11").Value = -36.5
$ws.Range("L11").Value = -19.6
$ws.Range("M11").Value = -24.87
$ws.Range("N11").Value = -13.85
$ws.Range("C12").Value = -48.4
$ws.Range("D12").Value = -34.9
$ws.Range("E12").Value = -50
$ws.Range("F12").Value = 89.3
$ws.Range("G12").Value = -39.1
$ws.Range("L12").Value = -19.6
$ws.Range("M12").Value = -24.87
$ws.Range("N12").Value = -13.85
$ws.Range("C13").Value = -32.4
$ws.Range("D13").Value = -45.2
$ws.Range("E13").Value = -70.3
$ws.Range("F13").Value = 32.8
$ws.Range("G13").Value = -79.8
$ws.Range("C14").Value = -57.5
$ws.Range("D14").Value = -40.5
$ws.Range("E14").Value = -62.7
$ws.Range("F14").Value = 45.6
$ws.Range("G14").Value = -62.5
$ws.Range("L14").Value = -37.1
$ws.Range("M14").Value = -44.48
$ws.Range("N14").Value = -35.07
$ws.Range("C15").Value = -34.5
$ws.Range("D15").Value = -40.3
$ws.Range("E15").Value = -74.7
$ws.Range("F15").Value = 40.8
$ws.Range("G15").Value = -74.90000000000001
$ws.Range("L15").Value = -26.98
$ws.Range("M15").Value = -48.38
$ws.Range("N15").Value = -50.92
$ws.Range("C16").Value = -29.9
$ws.Range("D16").Value = -37.9
$ws.Range("E16").Value = -76.5
$ws.Range("F16").Value = 45.1
$ws.Range("G16").Value = -74.90000000000001
$ws.Range("L16").Value = -9.210000000000001
$ws.Range("M16").Value = -29.29
$ws.Range("N16").Value = -21.35
$ws.Range("C17").Value = -54.3
$ws.Range("D17").Value = -51.7
$ws.Range("E17").Value = -81.8
$ws.Range("F17").Value = 4.8
$ws.Range("G17").Value = -97.5
$ws.Range("C18").Value = -41.2
$ws.Range("D18").Value = -49.8
$ws.Range("E18").Value = -80.59999999999999
$ws.Range("F18").Value = 5.6
$ws.Range("G18").Value = -97.2
$ws.Range("L18").Value = -5.68
$ws.Range("M18").Value = -9.220000000000001
$ws.Range("N18").Value = -11.55
$ws.Range("C19").Value = -33.8
$ws.Range("D19").Value = -36.6
$ws.Range("E19").Value = -67
$ws.Range("F19").Value = 28.9
$ws.Range("G19").Value = -79.5
$ws.Range("C20").Value = -43.1
$ws.Range("D20").Value = -45.4
$ws.Range("E20").Value = -84.8
$ws.Range("F20").Value = 26
$ws.Range("G20").Value = -87.2
$ws.Range("C21").Value = -30
$ws.Range("D21").Value = -36.3
$ws.Range("E21").Value = -75
$ws.Range("F21").Value = 48.2
$ws.Range("G21").Value = -71.7
$ws.Range("C22").Value = -37
$ws.Range("D22").Value = -32.8
$ws.Range("E22").Value = -75.3
$ws.Range("F22").Value = 32.4
$ws.Range("G22").Value = -77.3
$ws.Range("C23").Value = -33.4
$ws.Range("D23").Value = -50.5
$ws.Range("E23").Value = -76.2
$ws.Range("F23").Value = 13.1
$ws.Range("G23").Value = -91.90000000000001
$ws.Range("C24").Value = -38.7
$ws.Range("D24").Value = -53
$ws.Range("E24").Value = -81.3
$ws.Range("F24").Value = 20.8
$ws.Range("G24").Value = -84.8
$ws.Range("C25").Value = -50.1
$ws.Range("D25").Value = -42.1
$ws.Range("E25").Value = -57.4
$ws.Range("F25").Value = 33.7
$ws.Range("G25").Value = -72.3
$ws.Range("C26").Value = -33.2
$ws.Range("D26").Value = -44.1
$ws.Range("E26").Value = -77.8
$ws.Range("F26").Value = 6.3
$ws.Range("G26").Value = -95.7
$ws.Range("C27").Value = -56.9
$ws.Range("D27").Value = -51.3
$ws.Range("E27").Value = -75.59999999999999
$ws.Range("F27").Value = 40.4
$ws.Range("G27").Value = -76.59999999999999
$ws.Range("C28").Value = -27.3
$ws.Range("D28").Value = -41.9
$ws.Range("E28").Value = -73.09999999999999
$ws.Range("F28").Value = 39.8
$ws.Range("G28").Value = -71.09999999999999
$ws.Range("L28").Value = -46.96
$ws.Range("M28").Value = -64.04000000000001
$ws.Range("N28").Value = -59.66
$ws.Range("C29").Value = -54.6
$ws.Range("D29").Value = -44.7
$ws.Range("E29").Value = -70.5
$ws.Range("F29").Value = 42.9
$ws.Range("G29").Value = -67.40000000000001
$ws.Range("L29").Value = -56.55
$ws.Range("M29").Value = -68.41
$ws.Range("N29").Value = -64.04000000000001
$ws.Range("C30").Value = -23
$ws.Range("D30").Value = -33.5
$ws.Range("E30").Value = -75.7
$ws.Range("F30").Value = 54.3
$ws.Range("G30").Value = -69.59999999999999
$ws.Range("L30").Value = -46.96
$ws.Range("M30").Value = -64.04000000000001
$ws.Range("N30").Value = -59.66
$ws.Range("C31").Value = -28.9
$ws.Range("D31").Value = -41.1
$ws.Range("E31").Value = -67.5
$ws.Range("F31").Value = 50.6
$ws.Range("G31").Value = -66
$ws.Range("C32").Value = -27
$ws.Range("D32").Value = -42.4
$ws.Range("E32").Value = -72.90000000000001
$ws.Range("F32").Value = 39.8
$ws.Range("G32").Value = -70.90000000000001
$ws.Range("L32").Value = -46.96
$ws.Range("M32").Value = -64.04000000000001
$ws.Range("N32").Value = -59.66
$ws.Range("C33").Value = -38.1
$ws.Range("D33").Value = -50.8
$ws.Range("E33").Value = -76.2
$ws.Range("F33").Value = 19.4
$ws.Range("G33").Value = -87.59999999999999
$ws.Range("C34").Value = -56.6
$ws.Range("D34").Value = -38.7
$ws.Range("E34").Value = -60.7
$ws.Range("F34").Value = 51
$ws.Range("G34").Value = -56.9
$ws.Range("L34").Value = -37.1
$ws.Range("M34").Value = -44.48
$ws.Range("N34").Value = -35.07
$ws.Range("C35").Value = -59.4
$ws.Range("D35").Value = -52.4
$ws.Range("E35").Value = -76.90000000000001
$ws.Range("F35").Value = 35.9
$ws.Range("G35").Value = -79.90000000000001
$ws.Range("C36").Value = -62.3
$ws.Range("D36").Value = -47.1
$ws.Range("E36").Value = -64.40000000000001
$ws.Range("F36").Value = 62.2
$ws.Range("G36").Value = -65.8
$ws.Range("L36").Value = -78.22
$ws.Range("M36").Value = -77.61
$ws.Range("N36").Value = -70.23999999999999
$ws.Range("C37").Value = -24.8
$ws.Range("D37").Value = -33
$ws.Range("E37").Value = -69.2
$ws.Range("F37").Value = 60.6
$ws.Range("G37").Value = -61.7
$ws.Range("L37").Value = -46.96
$ws.Range("M37").Value = -64.04000000000001
$ws.Range("N37").Value = -59.66
$ws.Range("C38").Value = -43.5
$ws.Range("D38").Value = -31.5
$ws.Range("E38").Value = -50.6
$ws.Range("F38").Value = 61.3
$ws.Range("G38").Value = -38.7
$ws.Range("L38").Value = -22.57
$ws.Range("M38").Value = -42.3
$ws.Range("N38").Value = -33.04
$ws.Range("C39").Value = -21.6
$ws.Range("D39").Value = -30.9
$ws.Range("E39").Value = -74.59999999999999
$ws.Range("F39").Value = 55.2
$ws.Range("G39").Value = -68.5
$ws.Range("L39").Value = -46.96
$ws.Range("M39").Value = -64.04000000000001
$ws.Range("N39").Value = -59.66
$ws.Range("C40").Value = -26.4
$ws.Range("D40").Value = -35.5
$ws.Range("E40").Value = -67.59999999999999
$ws.Range("F40").Value = 59.9
$ws.Range("G40").Value = -60.2
$ws.Range("L40").Value = -46.96
$ws.Range("M40").Value = -64.04000000000001
$ws.Range("N40").Value = -59.66
$ws.Range("C41").Value = -53.5
$ws.Range("D41").Value = -38.2
$ws.Range("E41").Value = -57
$ws.Range("F41").Value = 76.2
$ws.Range("G41").Value = -44.4
$ws.Range("L41").Value = -37.1
$ws.Range("M41").Value = -44.48
$ws.Range("N41").Value = -35.07
$ws.Range("C42").Value = -35.9
$ws.Range("D42").Value = -59.1
$ws.Range("E42").Value = -85
$ws.Range("F42").Value = 18.7
$ws.Range("G42").Value = -87.09999999999999
$ws.Range("L42").Value = -6.29
$ws.Range("M42").Value = -4.1
$ws.Range("N42").Value = -16.53
